$d = $word.ActiveDocument

$d.Content.Find.Execute("(123) 456-7890", $true, $false, $false, $false, $false,
                         $true, 1, $false, "电话：（123） 456-7890", 2)

$d.Content.Find.Execute("首席动画师（2018 年 1 月至今）", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ABC 工作室：首席动画师（2018 年 1 月 - 演示）", 2)

$d.Content.Find.Execute("高级动画师（2015 年 6 月 - 2017 年 12 月）", $true, $false, $false, $false, $false,
                         $true, 1, $false, "XYZ 媒体：高级动画师 （2015 年 6 月 - 2017 年 12 月）", 2)

$d.Content.Find.Execute("初级动画师（2012 年 9 月 - 2015 年 5 月）", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MNO 娱乐： 初级动画师 （2012 年 9 月 - 2015 年 5 月）", 2)
